# Update crypto price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.711.97'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.300.67'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").Value = "'554.28"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = "'182.96"
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '3.290.51'
$ws.Range("E8").Value = '  +1.03%  '
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("E10").Value = '  -6.58%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = "'45.27"
$ws.Range("E12").Value = '  -3.87%  '
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").Value = '3.825.69'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = "'8.31"
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").Value = "'573.15"
$ws.Range("E16").Value = '  -9.38%  '
$ws.Range("D17").Value = '65.589.36'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '3.297.52'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("D20").Value = "'17.44"
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").Value = "'10.72"
$ws.Range("E21").Value = '  -5.00%  '
$ws.Range("D22").Value = "'0.880"
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("D23").Value = "'17.59"
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("D24").Value = "'4.95"
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").Value = "'97.38"
$ws.Range("E25").Value = '  -8.83%  '
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = "'9.18"
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("D30").Value = "'8.31"
$ws.Range("E30").Value = '  -3.99%  '
$ws.Range("D31").Value = "'30.17"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").Value = "'6.50"
$ws.Range("E32").Value = '  +4.49%  '
$ws.Range("D33").Value = "'556.35"
$ws.Range("E33").Value = '  +6.82%  '
$ws.Range("E34").Value = '  -7.49%  '
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("D36").Value = '3.755.76'
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("D39").Value = "'55.44"
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("D40").Value = "'32.93"
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("E41").Value = '  -4.36%  '
$ws.Range("E42").Value = '  -8.87%  '
$ws.Range("E43").Value = '  +2.78%  '
$ws.Range("D44").Value = '0.0₃0668'
$ws.Range("E44").Value = '  -8.58%  '
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = '  -6.27%  '
$ws.Range("D46").Value = "'0.328"
$ws.Range("E46").Value = '  -2.42%  '
$ws.Range("E47").Value = '  -2.69%  '
$ws.Range("D48").Value = "'3.02"
$ws.Range("E48").Value = '  -10.98%  '
$ws.Range("D49").Value = "'0.998"
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("E51").Value = '  -4.83%  '
